# Update column C ("Förändrad") from serial date 45177 to 45178 for all data rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row  # xlUp = -4162

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45177) {
        $cell.Value2 = 45178
    }
}
